$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "PERALTA REYES MARY CRUZ",
    "ESPINOZA GUZMAN MAYRA LOURDES",
    "RUIZ CARRASCO HILLARY SAMANTHA",
    "CHUNGA DE LA CRUZ ROSA LILIANA",
    "FIESTAS PERICHE VIVIANA LISSETH",
    "ROSILLO ALBERCA ROXANA",
    "PRADO ACARO VANESSA PAOLA",
    "PALMA CARMENES DE MENA MERCEDES EVERJISTA",
    "PAIVA PINDAY ALICIA",
    "JIMENEZ GUERRERO JUAN RICARDO",
    "GIRON SILUPU JUAN FRANCISCO",
    "PAIVA GARCIA DANIELA MILEYDI"
)

$totals = @(243, 172, 168, 156, 155, 136, 118, 114, 113, 109, 103, 97)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $totals[$i]
}
